$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 7260
$ws.Range("I82").Value = 6350
$ws.Range("K82").Value = 19050
$ws.Range("M82").Value = -18644

$ws.Range("H85").Value = 7260
$ws.Range("I85").Value = 6350
$ws.Range("K85").Value = 19050
$ws.Range("M85").Value = -17646

$ws.Range("H92").Value = 699.2857
$ws.Range("I92").Value = 526.36365
$ws.Range("J92").Value = 1333.3334
$ws.Range("K92").Value = 526.36365
$ws.Range("L92").Value = 1333.3334
$ws.Range("M92").Value = 721.63635
$ws.Range("N92").Value = -3829.3334

$ws.Range("H97").Value = 889.6667
$ws.Range("J97").Value = 889.6667
$ws.Range("L97").Value = 2669.0001
$ws.Range("N97").Value = -3661.0001

$ws.Range("H99").Value = 1240.091
$ws.Range("I99").Value = 392.625
$ws.Range("K99").Value = 1177.875
$ws.Range("M99").Value = 320.125

$ws.Range("H133").Value = 26890
$ws.Range("J133").Value = 26890
$ws.Range("L133").Value = 26890
$ws.Range("N133").Value = -37010

$ws.Range("H135").Value = 1296.7693
$ws.Range("I135").Value = 1014.0909
$ws.Range("J135").Value = 2851.5
$ws.Range("K135").Value = 9126.8181
$ws.Range("L135").Value = 25663.5
$ws.Range("M135").Value = -6591.8181
$ws.Range("N135").Value = -30733.5

$ws.Range("H137").Value = 3642.0344
$ws.Range("I137").Value = 3438.2917
$ws.Range("J137").Value = 4620
$ws.Range("K137").Value = 10314.8751
$ws.Range("L137").Value = 13860
$ws.Range("M137").Value = -7764.875100000001
$ws.Range("N137").Value = -18960

$ws.Range("H138").Value = 2754.72
$ws.Range("I138").Value = 1541.875
$ws.Range("J138").Value = 4910.8887
$ws.Range("K138").Value = 4625.625
$ws.Range("L138").Value = 14732.6661
$ws.Range("M138").Value = 514.375
$ws.Range("N138").Value = -25012.6661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2312
$ws.Range("I61").Value = 1131.6
$ws.Range("J61").Value = 4047.8823
$ws.Range("K61").Value = 1131.6
$ws.Range("L61").Value = 4047.8823
$ws.Range("M61").Value = -919.5999999999999
$ws.Range("N61").Value = -4471.8823

$ws.Range("H122").Value = 2480.75
$ws.Range("I122").Value = 1634.6666
$ws.Range("J122").Value = 3568.5715
$ws.Range("K122").Value = 4903.9998
$ws.Range("L122").Value = 10705.7145
$ws.Range("M122").Value = -2453.9998
$ws.Range("N122").Value = -15605.7145

$ws.Range("H132").Value = 2075.8333
$ws.Range("I132").Value = 1445.9778
$ws.Range("J132").Value = 5225.1113
$ws.Range("K132").Value = 4337.9334
$ws.Range("L132").Value = 15675.3339
$ws.Range("M132").Value = -1807.9334
$ws.Range("N132").Value = -20735.3339

$ws.Range("H136").Value = 2312
$ws.Range("I136").Value = 1131.6
$ws.Range("J136").Value = 4047.8823
$ws.Range("K136").Value = 3394.8
$ws.Range("L136").Value = 12143.6469
$ws.Range("M136").Value = -844.7999999999997
$ws.Range("N136").Value = -17243.6469

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2708.3076
$ws.Range("I31").Value = 2120
$ws.Range("J31").Value = 5944
$ws.Range("K31").Value = 2120
$ws.Range("L31").Value = 5944
$ws.Range("M31").Value = -1825
$ws.Range("N31").Value = -6534

$ws.Range("H34").Value = 2708.3076
$ws.Range("I34").Value = 2120
$ws.Range("J34").Value = 5944
$ws.Range("K34").Value = 2120
$ws.Range("L34").Value = 5944
$ws.Range("M34").Value = -1918
$ws.Range("N34").Value = -6348

$ws.Range("H99").Value = 3542.7144
$ws.Range("I99").Value = 1500
$ws.Range("J99").Value = 4359.8
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 4359.8
$ws.Range("M99").Value = -2
$ws.Range("N99").Value = -7355.8

$ws.Range("H126").Value = 3542.7144
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 4359.8
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 13079.4
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -18019.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1167.5493
$ws.Range("I131").Value = 2948
$ws.Range("J131").Value = 1032.6666
$ws.Range("K131").Value = 8844
$ws.Range("L131").Value = 3097.9998
$ws.Range("M131").Value = -3804
$ws.Range("N131").Value = -13177.9998

$ws.Range("H138").Value = 1741.2727
$ws.Range("I138").Value = 1021.8333
$ws.Range("K138").Value = 3065.4999
$ws.Range("M138").Value = 2074.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3451.5925
$ws.Range("I122").Value = 2899.611
$ws.Range("K122").Value = 8698.832999999999
$ws.Range("M122").Value = -6248.832999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 680.5217
$ws.Range("I16").Value = 680.5217
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 680.5217
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -510.5217
$ws.Range("N16").ClearContents()

$ws.Range("H93").Value = 2558
$ws.Range("I93").Value = 1541.1428
$ws.Range("K93").Value = 1541.1428
$ws.Range("M93").Value = -293.1428000000001

$ws.Range("H122").Value = 3864.2856
$ws.Range("J122").Value = 4942.857
$ws.Range("L122").Value = 14828.571
$ws.Range("N122").Value = -19728.571

$ws.Range("H136").Value = 1766.921
$ws.Range("I136").Value = 1283.4572
$ws.Range("J136").Value = 7407.3335
$ws.Range("K136").Value = 3850.3716
$ws.Range("L136").Value = 22222.0005
$ws.Range("M136").Value = -1300.3716
$ws.Range("N136").Value = -27322.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 20153
$ws.Range("J6").Value = 20153
$ws.Range("L6").Value = 20153
$ws.Range("N6").Value = -20383

$ws.Range("H132").Value = 3456.224
$ws.Range("I132").Value = 1242.5652
$ws.Range("J132").Value = 11941.917
$ws.Range("K132").Value = 3727.6956
$ws.Range("L132").Value = 35825.751
$ws.Range("M132").Value = -1197.6956
$ws.Range("N132").Value = -40885.751

$ws.Range("H136").Value = 912.4583
$ws.Range("I136").Value = 574.9545000000001
$ws.Range("J136").Value = 4625
$ws.Range("K136").Value = 1724.8635
$ws.Range("L136").Value = 13875
$ws.Range("M136").Value = 825.1364999999998
$ws.Range("N136").Value = -18975
